$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 120.42857
$ws.Cells.Item(12, 10).Value = 107.5
$ws.Cells.Item(12, 12).Value = 107.5
$ws.Cells.Item(12, 14).Value = -447.5

$ws.Cells.Item(28, 8).Value = 967.75
$ws.Cells.Item(28, 9).Value = 967.75
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 967.75
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = -482.75
$ws.Cells.Item(28, 14).ClearContents()

$ws.Cells.Item(62, 8).Value = 6394.8125
$ws.Cells.Item(62, 9).Value = 4434.5
$ws.Cells.Item(62, 10).Value = 7048.25
$ws.Cells.Item(62, 11).Value = 4434.5
$ws.Cells.Item(62, 12).Value = 7048.25
$ws.Cells.Item(62, 13).Value = -3810.5
$ws.Cells.Item(62, 14).Value = -8296.25

$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()

$ws.Cells.Item(65, 8).Value = 6394.8125
$ws.Cells.Item(65, 9).Value = 4434.5
$ws.Cells.Item(65, 10).Value = 7048.25
$ws.Cells.Item(65, 11).Value = 22172.5
$ws.Cells.Item(65, 12).Value = 35241.25
$ws.Cells.Item(65, 13).Value = -19052.5
$ws.Cells.Item(65, 14).Value = -41481.25

$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()

$ws.Cells.Item(76, 8).Value = 6445.636
$ws.Cells.Item(76, 10).Value = 7166.8887
$ws.Cells.Item(76, 12).Value = 7166.8887
$ws.Cells.Item(76, 14).Value = -7796.8887

$ws.Cells.Item(79, 8).Value = 6445.636
$ws.Cells.Item(79, 10).Value = 7166.8887
$ws.Cells.Item(79, 12).Value = 7166.8887
$ws.Cells.Item(79, 14).Value = -9350.8887

$ws.Cells.Item(107, 8).Value = 1336
$ws.Cells.Item(107, 9).Value = 1311.7646
$ws.Cells.Item(107, 11).Value = 1311.7646
$ws.Cells.Item(107, 13).Value = 608.2354

$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 11216.5
$ws.Cells.Item(28, 9).Value = 11216.5
$ws.Cells.Item(28, 11).Value = 11216.5
$ws.Cells.Item(28, 13).Value = -11024.5

$ws.Cells.Item(45, 8).Value = 3807.65
$ws.Cells.Item(45, 9).Value = 2381.9092
$ws.Cells.Item(45, 10).Value = 5550.222
$ws.Cells.Item(45, 11).Value = 2381.9092
$ws.Cells.Item(45, 12).Value = 5550.222
$ws.Cells.Item(45, 13).Value = -2004.9092
$ws.Cells.Item(45, 14).Value = -6304.222

$ws.Cells.Item(61, 8).Value = 3950.5518
$ws.Cells.Item(61, 9).Value = 2472.5
$ws.Cells.Item(61, 10).Value = 11045.2
$ws.Cells.Item(61, 11).Value = 2472.5
$ws.Cells.Item(61, 12).Value = 11045.2
$ws.Cells.Item(61, 13).Value = -2260.5
$ws.Cells.Item(61, 14).Value = -11469.2

$ws.Cells.Item(99, 8).Value = 11216.5
$ws.Cells.Item(99, 9).Value = 11216.5
$ws.Cells.Item(99, 11).Value = 11216.5
$ws.Cells.Item(99, 13).Value = -8221.5

$ws.Cells.Item(110, 8).Value = 1251.2727
$ws.Cells.Item(110, 9).Value = 1254.1666
$ws.Cells.Item(110, 11).Value = 1254.1666
$ws.Cells.Item(110, 13).Value = 790.8334

$ws.Cells.Item(136, 8).Value = 3950.5518
$ws.Cells.Item(136, 9).Value = 2472.5
$ws.Cells.Item(136, 10).Value = 11045.2
$ws.Cells.Item(136, 11).Value = 7417.5
$ws.Cells.Item(136, 12).Value = 33135.60000000001
$ws.Cells.Item(136, 13).Value = -4867.5
$ws.Cells.Item(136, 14).Value = -38235.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 43742.453
$ws.Cells.Item(86, 9).Value = 29000
$ws.Cells.Item(86, 10).Value = 47018.555
$ws.Cells.Item(86, 11).Value = 29000
$ws.Cells.Item(86, 12).Value = 47018.555
$ws.Cells.Item(86, 13).Value = -27877
$ws.Cells.Item(86, 14).Value = -49264.555

$ws.Cells.Item(89, 8).Value = 43742.453
$ws.Cells.Item(89, 9).Value = 29000
$ws.Cells.Item(89, 10).Value = 47018.555
$ws.Cells.Item(89, 11).Value = 145000
$ws.Cells.Item(89, 12).Value = 235092.775
$ws.Cells.Item(89, 13).Value = -139384
$ws.Cells.Item(89, 14).Value = -246324.775

$ws.Cells.Item(97, 8).Value = 58247.5
$ws.Cells.Item(97, 10).Value = 58247.5
$ws.Cells.Item(97, 12).Value = 58247.5
$ws.Cells.Item(97, 14).Value = -60229.5

$ws.Cells.Item(99, 8).Value = 9549.4
$ws.Cells.Item(99, 9).Value = 10699.5
$ws.Cells.Item(99, 10).Value = 4949
$ws.Cells.Item(99, 11).Value = 10699.5
$ws.Cells.Item(99, 12).Value = 4949
$ws.Cells.Item(99, 13).Value = -9201.5
$ws.Cells.Item(99, 14).Value = -7945

$ws.Cells.Item(109, 8).Value = 25000
$ws.Cells.Item(109, 9).Value = 25000
$ws.Cells.Item(109, 11).Value = 25000
$ws.Cells.Item(109, 13).Value = -23960

$ws.Cells.Item(126, 8).Value = 9549.4
$ws.Cells.Item(126, 9).Value = 10699.5
$ws.Cells.Item(126, 10).Value = 4949
$ws.Cells.Item(126, 11).Value = 32098.5
$ws.Cells.Item(126, 12).Value = 14847
$ws.Cells.Item(126, 13).Value = -29628.5
$ws.Cells.Item(126, 14).Value = -19787

$ws.Cells.Item(132, 8).Value = 2721.077
$ws.Cells.Item(132, 9).Value = 2531.1667
$ws.Cells.Item(132, 11).Value = 7593.500100000001
$ws.Cells.Item(132, 13).Value = -5063.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 300
$ws.Cells.Item(87, 9).Value = 300
$ws.Cells.Item(87, 11).Value = 900
$ws.Cells.Item(87, 13).Value = 348

$ws.Cells.Item(90, 8).Value = 300
$ws.Cells.Item(90, 9).Value = 300
$ws.Cells.Item(90, 11).Value = 2700
$ws.Cells.Item(90, 13).Value = 3540

$ws.Cells.Item(117, 8).Value = 1042.7142
$ws.Cells.Item(117, 10).Value = 1042.7142
$ws.Cells.Item(117, 12).Value = 3128.1426
$ws.Cells.Item(117, 14).Value = -10012.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 12625.75
$ws.Cells.Item(10, 9).Value = 11126.5
$ws.Cells.Item(10, 10).Value = 13375.375
$ws.Cells.Item(10, 11).Value = 11126.5
$ws.Cells.Item(10, 12).Value = 13375.375
$ws.Cells.Item(10, 13).Value = -10957.5
$ws.Cells.Item(10, 14).Value = -13713.375

$ws.Cells.Item(12, 8).Value = 139909.1
$ws.Cells.Item(12, 9).Value = 3500
$ws.Cells.Item(12, 10).Value = 303600
$ws.Cells.Item(12, 11).Value = 3500
$ws.Cells.Item(12, 12).Value = 303600
$ws.Cells.Item(12, 13).Value = -3360
$ws.Cells.Item(12, 14).Value = -303880

$ws.Cells.Item(14, 8).Value = 2755.6086
$ws.Cells.Item(14, 10).Value = 2802
$ws.Cells.Item(14, 12).Value = 2802
$ws.Cells.Item(14, 14).Value = -3138

$ws.Cells.Item(63, 8).Value = 56666.668
$ws.Cells.Item(63, 10).Value = 56666.668
$ws.Cells.Item(63, 12).Value = 56666.668
$ws.Cells.Item(63, 14).Value = -58038.668

$ws.Cells.Item(66, 8).Value = 56666.668
$ws.Cells.Item(66, 10).Value = 56666.668
$ws.Cells.Item(66, 12).Value = 170000.004
$ws.Cells.Item(66, 14).Value = -176864.004

$ws.Cells.Item(88, 8).Value = 255555
$ws.Cells.Item(88, 10).Value = 255555
$ws.Cells.Item(88, 12).Value = 255555
$ws.Cells.Item(88, 14).Value = -256457

$ws.Cells.Item(91, 8).Value = 255555
$ws.Cells.Item(91, 10).Value = 255555
$ws.Cells.Item(91, 12).Value = 255555
$ws.Cells.Item(91, 14).Value = -258675

$ws.Cells.Item(101, 8).Value = 42740.715
$ws.Cells.Item(101, 10).Value = 42740.715
$ws.Cells.Item(101, 12).Value = 42740.715
$ws.Cells.Item(101, 14).Value = -49230.715

$ws.Cells.Item(102, 8).Value = 2094.087
$ws.Cells.Item(102, 9).Value = 1021.6923
$ws.Cells.Item(102, 11).Value = 1021.6923
$ws.Cells.Item(102, 13).Value = 600.3077

$ws.Cells.Item(122, 8).Value = 131197.62
$ws.Cells.Item(122, 9).Value = 6931.1665
$ws.Cells.Item(122, 10).Value = 503997
$ws.Cells.Item(122, 11).Value = 20793.4995
$ws.Cells.Item(122, 12).Value = 1511991
$ws.Cells.Item(122, 13).Value = -18343.4995
$ws.Cells.Item(122, 14).Value = -1516891

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 7366.4614
$ws.Cells.Item(100, 9).Value = 5796.909
$ws.Cells.Item(100, 11).Value = 5796.909
$ws.Cells.Item(100, 13).Value = -5255.909

$ws.Cells.Item(122, 8).Value = 27685.54
$ws.Cells.Item(122, 9).Value = 26102.445
$ws.Cells.Item(122, 11).Value = 78307.33499999999
$ws.Cells.Item(122, 13).Value = -75857.33499999999

$ws.Cells.Item(132, 8).Value = 16167.667
$ws.Cells.Item(132, 10).Value = 8999
$ws.Cells.Item(132, 12).Value = 26997
$ws.Cells.Item(132, 14).Value = -32057

$ws.Cells.Item(136, 8).Value = 16642.285
$ws.Cells.Item(136, 9).Value = 12399.6
$ws.Cells.Item(136, 11).Value = 37198.8
$ws.Cells.Item(136, 13).Value = -34648.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 14).ClearContents()

$ws.Cells.Item(15, 8).Value = 40007
$ws.Cells.Item(15, 10).Value = 40007
$ws.Cells.Item(15, 12).Value = 40007
$ws.Cells.Item(15, 14).Value = -40583

$ws.Cells.Item(58, 8).Value = 15995
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(62, 8).Value = 4792.4526
$ws.Cells.Item(62, 9).Value = 3885.7144
$ws.Cells.Item(62, 11).Value = 3885.7144
$ws.Cells.Item(62, 13).Value = -3261.7144

$ws.Cells.Item(65, 8).Value = 4792.4526
$ws.Cells.Item(65, 9).Value = 3885.7144
$ws.Cells.Item(65, 11).Value = 19428.572
$ws.Cells.Item(65, 13).Value = -16308.572

$ws.Cells.Item(81, 8).Value = 1929.5
$ws.Cells.Item(81, 10).Value = 1700
$ws.Cells.Item(81, 12).Value = 3400
$ws.Cells.Item(81, 14).Value = -5522

$ws.Cells.Item(84, 8).Value = 1929.5
$ws.Cells.Item(84, 10).Value = 1700
$ws.Cells.Item(84, 12).Value = 17000
$ws.Cells.Item(84, 14).Value = -27608

$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 1392.091
$ws.Cells.Item(126, 9).Value = 1392.091
$ws.Cells.Item(126, 11).Value = 4176.272999999999
$ws.Cells.Item(126, 13).Value = -1706.272999999999

$ws.Cells.Item(136, 8).Value = 5034.8
$ws.Cells.Item(136, 9).Value = 5034.8
$ws.Cells.Item(136, 11).Value = 15104.4
$ws.Cells.Item(136, 13).Value = -12554.4
